$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values must be stored as literal text (not auto-converted
# to numbers) to preserve formatting like trailing zeros and thousand-dot
# separators, matching the original inlineStr cell content.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.626.58'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.614.20'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '627.93'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.14'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.612.15'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.496'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.39'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.440'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000228'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.34'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.233.91'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.818.43'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.617.34'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.01'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.18'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '464.01'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.645'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.74'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.70'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.761.24'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.18'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.72'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.177'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.56'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '26.56'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.96'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.615.35'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.46'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.41'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '179.95'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0926'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.68'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '32.49'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.914'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '46.03'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.81'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.268'

# Restore the default (Normal) cell style so only the value changed.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"

# Coin name / link / volume columns are plain text already, no coercion risk.
$ws.Range("E2").Value = '  +4.77%  '
$ws.Range("E3").Value = '  +4.15%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E5").Value = '  +4.42%  '
$ws.Range("E6").Value = '  +6.86%  '
$ws.Range("E7").Value = '  +4.50%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("E9").Value = '  +4.22%  '
$ws.Range("E10").Value = '  +8.68%  '
$ws.Range("E11").Value = '  +8.10%  '
$ws.Range("E12").Value = '  +5.34%  '
$ws.Range("E13").Value = '  +5.90%  '
$ws.Range("E14").Value = '  +7.33%  '
$ws.Range("E15").Value = '  +4.29%  '
$ws.Range("E16").Value = '  +4.97%  '
$ws.Range("E17").Value = '  +3.09%  '
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("E19").Value = '  +6.47%  '
$ws.Range("E20").Value = '  +7.08%  '
$ws.Range("E21").Value = '  +15.23%  '
$ws.Range("E22").Value = '  +5.73%  '
$ws.Range("E23").Value = '  +4.74%  '
$ws.Range("E24").Value = '  +2.64%  '
$ws.Range("E25").Value = '  +14.02%  '
$ws.Range("E26").Value = '  +7.30%  '
$ws.Range("E27").Value = '  +3.92%  '
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("E29").Value = '  +14.17%  '
$ws.Range("E30").Value = '  +5.76%  '
$ws.Range("E31").Value = '  +10.40%  '
$ws.Range("E32").Value = '  +13.40%  '
$ws.Range("E33").Value = '  +8.25%  '
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("E35").Value = '  +5.01%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("E36").Value = '  +7.19%  '
$ws.Range("E37").Value = '  +4.38%  '
$ws.Range("E38").Value = '  +7.52%  '
$ws.Range("E39").Value = '  +13.06%  '
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("E41").Value = '  +5.07%  '
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("E42").Value = '  +9.14%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("E44").Value = '  +4.88%  '
$ws.Range("E45").Value = '  +21.03%  '
$ws.Range("E46").Value = '  +4.09%  '
$ws.Range("E47").Value = '  +15.16%  '
$ws.Range("E48").Value = '  +2.21%  '
$ws.Range("E49").Value = '  +13.62%  '
$ws.Range("E50").Value = '  +5.04%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("E51").Value = '  +10.07%  '
